$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index map for the columns touched by this edit
# C=3 block_total, F=6 trial_total, H=8 category, I=9 cond_cat,
# K=11 correct_answer, L=12 stimulus, M=13 conceptual, N=14 perceptual,
# O=15 typicality, P=16 n, Q=17 p_typicality, R=18 p_conceptual, S=19 p_perceptual

$data = @(
  @{ row=2; C=3; F=135; H='bedrooms'; I='distractor'; K='f'; L='stimuli/img_fea1z.png'; M=79.45945945945945; N=56.24324324324324; O=67.85135135135135; P=37; Q=7; R=7; S=7 }
  @{ row=3; C=3; F=136; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_bbs77.png'; M=31.64444444444445; N=21.26666666666667; O=26.45555555555556; P=45; Q=2; R=2; S=2 }
  @{ row=4; C=3; F=137; H='bedrooms'; I='distractor'; K='f'; L='stimuli/img_ys3qz.png'; M=46.79545454545455; N=31.20454545454545; O=39; P=44; Q=2; R=2; S=2 }
  @{ row=5; C=3; F=138; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_6zz63.png'; M=87.66666666666667; N=70.6; O=79.13333333333333; P=45; Q=9; R=10; S=10 }
  @{ row=6; C=3; F=139; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_xy930.png'; M=70.5952380952381; N=49.47619047619047; O=60.03571428571429; P=42; Q=6; R=6; S=6 }
  @{ row=7; C=3; F=140; H='kitchens'; I='distractor'; K='f'; L='stimuli/img_cxpff.png'; M=74.92307692307692; N=53.28205128205128; O=64.1025641025641; P=39; Q=6; R=6; S=6 }
  @{ row=8; C=3; F=141; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_wz6x5.png'; M=68.3695652173913; N=48.47826086956522; O=58.42391304347826; P=46; Q=5; R=5; S=5 }
  @{ row=9; C=3; F=142; H='living_rooms'; I='target'; K='j'; L='stimuli/img_eh0no.png'; M=53.66666666666666; N=36.02564102564103; O=44.84615384615385; P=39; Q=3; R=3; S=3 }
  @{ row=10; C=3; F=143; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_16kib.png'; M=80.97727272727273; N=61.11363636363637; O=71.04545454545455; P=44; Q=8; R=8; S=8 }
  @{ row=11; C=3; F=144; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_xu1p3.png'; M=75.27659574468085; N=56.68085106382978; O=65.97872340425532; P=47; Q=7; R=7; S=7 }
  @{ row=12; C=3; F=145; H='bedrooms'; I='distractor'; K='f'; L='stimuli/img_kugyw.png'; M=74.25; N=54.10714285714285; O=64.17857142857143; P=28; Q=6; R=6; S=6 }
  @{ row=13; C=3; F=146; H='bedrooms'; I='target'; K='j'; L='stimuli/img_4o8l0.png'; M=46.02173913043478; N=31.45652173913043; O=38.73913043478261; P=46; Q=3; R=3; S=3 }
  @{ row=14; C=3; F=147; H='bedrooms'; I='target'; K='j'; L='stimuli/img_kost0.png'; M=63.09090909090909; N=42.77272727272727; O=52.93181818181819; P=44; Q=5; R=5; S=5 }
  @{ row=15; C=3; F=148; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_w8yhd.png'; M=55.74418604651163; N=38.90697674418605; O=47.32558139534883; P=43; Q=4; R=4; S=4 }
  @{ row=16; C=3; F=149; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_abobq.png'; M=75.1842105263158; N=54.13157894736842; O=64.65789473684211; P=38; Q=6; R=6; S=6 }
  @{ row=17; C=3; F=150; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_0kqc0.png'; M=43.74468085106383; N=27.14893617021277; O=35.4468085106383; P=47; Q=2; R=2; S=2 }
  @{ row=18; C=3; F=151; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_wgkqa.png'; M=87.25581395348837; N=71.13953488372093; O=79.19767441860465; P=43; Q=10; R=10; S=10 }
  @{ row=19; C=3; F=152; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_bj99b.png'; M=82.79069767441861; N=65.46511627906976; O=74.12790697674419; P=43; Q=8; R=8; S=8 }
  @{ row=20; C=3; F=153; H='bedrooms'; I='target'; K='j'; L='stimuli/img_amsgw.png'; M=86.08510638297872; N=65.95744680851064; O=76.02127659574468; P=47; Q=9; R=9; S=9 }
  @{ row=21; C=3; F=154; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_di6f0.png'; M=94.04347826086956; N=83.34782608695652; O=88.69565217391303; P=46; Q=10; R=10; S=10 }
  @{ row=22; C=3; F=155; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_cehin.png'; M=78.86363636363636; N=60.02272727272727; O=69.44318181818181; P=44; Q=7; R=7; S=7 }
  @{ row=23; C=3; F=156; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_xbtev.png'; M=13.68181818181818; N=8.568181818181818; O=11.125; P=44; Q=1; R=1; S=1 }
  @{ row=24; C=3; F=157; H='bedrooms'; I='target'; K='j'; L='stimuli/img_pey7u.png'; M=30.34883720930232; N=20.34883720930232; O=25.34883720930232; P=43; Q=1; R=2; S=2 }
  @{ row=25; C=3; F=158; H='living_rooms'; I='distractor'; K='f'; L='stimuli/img_6a0hu.png'; M=61.275; N=42.025; O=51.65; P=40; Q=4; R=4; S=4 }
  @{ row=26; C=3; F=159; H='bedrooms'; I='distractor'; K='f'; L='stimuli/img_pt3d7.png'; M=65.08571428571429; N=44.65714285714286; O=54.87142857142857; P=35; Q=4; R=4; S=4 }
  @{ row=27; C=3; F=160; H='bedrooms'; I='distractor'; K='f'; L='stimuli/img_twj5p.png'; M=67.71739130434783; N=42.08695652173913; O=54.90217391304348; P=46; Q=4; R=4; S=4 }
)

foreach ($d in $data) {
  $r = $d.row
  $ws.Cells.Item($r, 3).Value  = $d.C
  $ws.Cells.Item($r, 6).Value  = $d.F
  $ws.Cells.Item($r, 8).Value  = $d.H
  $ws.Cells.Item($r, 9).Value  = $d.I
  $ws.Cells.Item($r, 11).Value = $d.K
  $ws.Cells.Item($r, 12).Value = $d.L
  $ws.Cells.Item($r, 13).Value = $d.M
  $ws.Cells.Item($r, 14).Value = $d.N
  $ws.Cells.Item($r, 15).Value = $d.O
  $ws.Cells.Item($r, 16).Value = $d.P
  $ws.Cells.Item($r, 17).Value = $d.Q
  $ws.Cells.Item($r, 18).Value = $d.R
  $ws.Cells.Item($r, 19).Value = $d.S
}
